# Adds season-record columns (Wins / Losses / Ties) to the team stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: AD1 = "Wins", AE1 = "Losses", AF1 = "Ties"
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header formatting used by the rest of row 1 (bold, centered,
# bordered header style) by copying the format from the adjacent header cell.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# Every data row (2-45) gets the same season record: 90 wins, 72 losses, 0 ties.
$lastRow = 45
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 90  # column AD
    $ws.Cells.Item($r, 31).Value = 72  # column AE
    $ws.Cells.Item($r, 32).Value = 0   # column AF
}
